$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row (account 005273382 / MVFC / 185737.83) right after the
# ALPHASITIO row (row 2), pushing VIRGILIO and everything below down by one.
$ws.Rows.Item(3).Insert()
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "005273382"
$ws.Range("B3").Value = "MVFC"
$ws.Range("C3").Value = 185737.83

# The former "004267119 / ANA / 14593.13" row is now row 7 after the
# insertion above - replace it with "005305965 / SIDMAR / 16177.75".
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "005305965"
$ws.Range("B7").Value = "SIDMAR"
$ws.Range("C7").Value = 16177.75

# Remove the four rows that followed LUIZ (ISABEL, ISABEL, VIVIANE,
# PATRICIA), now at rows 9-12 after the earlier insertion.
$ws.Range("A9:A12").EntireRow.Delete()
